$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2 = @{ "B" = 16.77185617551796; "C" = 11.19200448434655; "D" = 4.564713903783011; "F" = 22.73933556382128; "G" = 3.619270867593013; "I" = 20.70275037715229; "L" = 10.65855822886157; "O" = 20.3287034907991 }
    3 = @{ "B" = 16.06047371218243; "C" = 10.92735726418366; "D" = 4.517223184870461; "F" = 22.78407899319903; "G" = 3.621419518540736; "I" = 20.85948178960559; "L" = 10.63083551722679; "O" = 20.42365773763223 }
    4 = @{ "B" = 15.60766123320989; "C" = 10.7610886995079; "D" = 4.487643017685982; "F" = 22.820267302597; "G" = 3.622808484356736; "I" = 20.96157876875253; "L" = 10.61598397171553; "O" = 20.48854784879627 }
    5 = @{ "B" = 15.41934539451825; "C" = 10.69245786991439; "D" = 4.475490345469015; "F" = 22.83719548963155; "G" = 3.623392078394731; "I" = 21.00465721136658; "L" = 10.61048190725094; "O" = 20.51663932561305 }
    6 = @{ "B" = 15.38785428706098; "C" = 10.68101116774348; "D" = 4.473466667363621; "F" = 22.84013780616319; "G" = 3.623490047132863; "I" = 21.01189929065758; "L" = 10.60960163180205; "O" = 20.52140319577031 }
    7 = @{ "B" = 15.6051365527722; "C" = 10.76016656227362; "D" = 4.487479511731964; "F" = 22.82048678549992; "G" = 3.622816283657953; "I" = 20.96215377680406; "L" = 10.61590753646051; "O" = 20.48892003783843 }
    8 = @{ "B" = 16.53002946146266; "C" = 11.10158049795088; "D" = 4.548430746672983; "F" = 22.75294828276099; "G" = 3.619997293125181; "I" = 20.75557365964006; "L" = 10.64855173596979; "O" = 20.3600714617722 }
    9 = @{ "B" = 18.20772120234537; "C" = 11.73811831582631; "D" = 4.664304686360434; "F" = 22.69006296868202; "G" = 3.615019615772198; "I" = 20.39705206144924; "L" = 10.72957498388209; "O" = 20.16003486282868 }
    10 = @{ "B" = 19.34739587456026; "C" = 12.18194724995736; "D" = 4.746801354076205; "F" = 22.68671672935099; "G" = 3.611694428238309; "I" = 20.16212851996731; "L" = 10.79915465103146; "O" = 20.04563918953278 }
    11 = @{ "B" = 19.8440524821722; "C" = 12.37798326763481; "D" = 4.783674719404924; "F" = 22.69456337185991; "G" = 3.610253019649851; "I" = 20.06146195018598; "L" = 10.83291314048009; "O" = 20.00076862163751 }
    12 = @{ "B" = 20.02888115820239; "C" = 12.45132165806471; "D" = 4.797536676528511; "F" = 22.69888438633383; "G" = 3.609717381255769; "I" = 20.02423619242999; "L" = 10.84599216176522; "O" = 19.98481581573049 }
    13 = @{ "B" = 19.98922088156001; "C" = 12.43556757763469; "D" = 4.794555866106098; "F" = 22.69789373201473; "G" = 3.609832288108359; "I" = 20.0322135971032; "L" = 10.84316234562427; "O" = 19.98820521946677 }
    14 = @{ "B" = 19.85932404767906; "C" = 12.38403510964624; "D" = 4.784817218431734; "F" = 22.69489180742745; "G" = 3.610208748395078; "I" = 20.05838142143293; "L" = 10.83398328328641; "O" = 19.99943531361058 }
    15 = @{ "B" = 19.77933296050298; "C" = 12.35235177884556; "D" = 4.778838632903069; "F" = 22.69322884982234; "G" = 3.610440666924306; "I" = 20.07452654278991; "L" = 10.82839908475334; "O" = 20.00644956844341 }
    16 = @{ "B" = 19.31448985257855; "C" = 12.16901320662301; "D" = 4.744377811671203; "F" = 22.68639270603088; "G" = 3.611790056683495; "I" = 20.16883233384398; "L" = 10.7969902105655; "O" = 20.04871651721334 }
    17 = @{ "B" = 19.02365844987236; "C" = 12.05499990635304; "D" = 4.723064445431404; "F" = 22.68460082807821; "G" = 3.6126360717826; "I" = 20.22827608331098; "L" = 10.77825619341188; "O" = 20.0764876716072 }
    18 = @{ "B" = 18.85433466166746; "C" = 11.98887340075571; "D" = 4.710744338112125; "F" = 22.6844519938284; "G" = 3.613129385427426; "I" = 20.26305007813264; "L" = 10.76767961851546; "O" = 20.09313519639147 }
    19 = @{ "B" = 18.79665695011903; "C" = 11.96639148635499; "D" = 4.706562663782811; "F" = 22.68455294956487; "G" = 3.613297566667472; "I" = 20.27492409534609; "L" = 10.76413292560245; "O" = 20.09888732743122 }
    20 = @{ "B" = 19.05483048887303; "C" = 12.06719404147556; "D" = 4.725339677515036; "F" = 22.68470028967712; "G" = 3.612545318140746; "I" = 20.22188778301895; "L" = 10.78022994580971; "O" = 20.07346154013379 }
    21 = @{ "B" = 19.89756675048416; "C" = 12.39919616696612; "D" = 4.787680496510363; "F" = 22.69573690537129; "G" = 3.610097896693592; "I" = 20.05067099815607; "L" = 10.83667144071622; "O" = 19.99610851151917 }
    22 = @{ "B" = 20.42939290678308; "C" = 12.61093360476079; "D" = 4.827830919106247; "F" = 22.71081671168256; "G" = 3.608557746839216; "I" = 19.94398609955871; "L" = 10.87527734641153; "O" = 19.951611612752 }
    23 = @{ "B" = 20.14731368190787; "C" = 12.4984216249335; "D" = 4.806458403814156; "F" = 22.70204818047834; "G" = 3.609374337430257; "I" = 20.00044769299154; "L" = 10.85451802172268; "O" = 19.97480364452612 }
    24 = @{ "B" = 19.04074420767801; "C" = 12.06168287159734; "D" = 4.724311252516001; "F" = 22.68465257795181; "G" = 3.612586326276469; "I" = 20.22477406791016; "L" = 10.77933700793666; "O" = 20.07482753184679 }
    25 = @{ "B" = 17.76955342246793; "C" = 11.56985146873197; "D" = 4.633391455367724; "F" = 22.69957569010398; "G" = 3.616307663566808; "I" = 20.48904664897264; "L" = 10.70586589240257; "O" = 20.20846424334196 }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
